$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "42.951.75"
$ws.Range("E2").Value2 = "  +0.44%  "
$ws.Range("D3").Value2 = "2.282.23"
$ws.Range("E3").Value2 = "  +0.20%  "
$ws.Range("E4").Value2 = "  -0.06%  "
$ws.Range("D5").Value2 = "'249.94"
$ws.Range("E5").Value2 = "  -0.59%  "
$ws.Range("D6").Value2 = "'0.644"
$ws.Range("E6").Value2 = "  +1.12%  "
$ws.Range("D7").Value2 = "'78.12"
$ws.Range("E7").Value2 = "  +8.59%  "
$ws.Range("E8").Value2 = "  +0.06%  "
$ws.Range("E9").Value2 = "  +1.41%  "
$ws.Range("D10").Value2 = "'40.83"
$ws.Range("E10").Value2 = "  +6.10%  "
$ws.Range("D11").Value2 = "'0.0975"
$ws.Range("E11").Value2 = "  +0.35%  "
$ws.Range("E12").Value2 = "  +0.25%  "
$ws.Range("D13").Value2 = "'0.106"
$ws.Range("E13").Value2 = "  +0.21%  "
$ws.Range("D14").Value2 = "2.622.71"
$ws.Range("E14").Value2 = "  +0.30%  "
$ws.Range("D15").Value2 = "'15.10"
$ws.Range("E15").Value2 = "  +0.93%  "
$ws.Range("D16").Value2 = "'0.871"
$ws.Range("E16").Value2 = "  -1.83%  "
$ws.Range("D17").Value2 = "2.285.16"
$ws.Range("E17").Value2 = "  +0.65%  "
$ws.Range("D18").Value2 = "42.870.41"
$ws.Range("E18").Value2 = "  +0.42%  "
$ws.Range("D19").Value2 = "0.0₃0997"
$ws.Range("E19").Value2 = "  -1.79%  "
$ws.Range("E20").Value2 = "  -1.30%  "
$ws.Range("D21").Value2 = "'72.22"
$ws.Range("E21").Value2 = "  -1.59%  "
$ws.Range("D22").Value2 = "'234.26"
$ws.Range("E22").Value2 = "  -0.93%  "
$ws.Range("D23").Value2 = "'2.18"
$ws.Range("E23").Value2 = "  +1.66%  "
$ws.Range("D24").Value2 = "'3.80"
$ws.Range("E24").Value2 = "  -1.39%  "
$ws.Range("E25").Value2 = "  +0.01%  "
$ws.Range("D26").Value2 = "'11.38"
$ws.Range("E26").Value2 = "  -2.36%  "
$ws.Range("D27").Value2 = "'2.34"
$ws.Range("E27").Value2 = "  -4.15%  "
$ws.Range("D28").Value2 = "'2.18"
$ws.Range("E28").Value2 = "  +2.24%  "
$ws.Range("D29").Value2 = "'167.80"
$ws.Range("E29").Value2 = "  -0.30%  "
$ws.Range("D30").Value2 = "'20.94"
$ws.Range("E30").Value2 = "  -0.39%  "
$ws.Range("D31").Value2 = "'6.46"
$ws.Range("E31").Value2 = "  +0.78%  "
$ws.Range("D32").Value2 = "'0.0856"
$ws.Range("E32").Value2 = "  +6.90%  "
$ws.Range("E33").Value2 = "  -4.39%  "
$ws.Range("D34").Value2 = "'30.28"
$ws.Range("E34").Value2 = "  -2.23%  "
$ws.Range("D35").Value2 = "'0.128"
$ws.Range("E35").Value2 = "  +1.03%  "
$ws.Range("D36").Value2 = "'4.58"
$ws.Range("E36").Value2 = "  -1.34%  "
$ws.Range("E37").Value2 = "  +0.24%  "
$ws.Range("E38").Value2 = "  -2.14%  "
$ws.Range("D39").Value2 = "'13.81"
$ws.Range("E39").Value2 = "  +3.78%  "
$ws.Range("E40").Value2 = "  -2.40%  "
$ws.Range("D41").Value2 = "'5.88"
$ws.Range("E41").Value2 = "  +0.52%  "
$ws.Range("D42").Value2 = "'112.81"
$ws.Range("E42").Value2 = "  +17.59%  "
$ws.Range("E43").Value2 = "  -1.09%  "
$ws.Range("D44").Value2 = "'61.24"
$ws.Range("E44").Value2 = "  -0.53%  "
$ws.Range("D45").Value2 = "'8.91"
$ws.Range("E45").Value2 = "  -2.63%  "
$ws.Range("E46").Value2 = "  -1.08%  "
$ws.Range("E47").Value2 = "  -0.06%  "
$ws.Range("D48").Value2 = "'4.58"
$ws.Range("E48").Value2 = "  -7.87%  "
$ws.Range("E49").Value2 = "  -2.60%  "
$ws.Range("E50").Value2 = "  -1.98%  "
$ws.Range("D51").Value2 = "'4.25"
$ws.Range("E51").Value2 = "  -0.15%  "
